# Burndown.xlsx — "fixing dates and adding new task time"
#
# 1) Tempo!N4: total sprint hours drops from 150 to 50 (all the dependent
#    Tempo!N5:N8 / O4:O8 formulas recompute automatically, as do the
#    Sprint 4 rows 14/15 that ultimately reference Tempo!N4).
# 2) Sprint 4!J2: the date header that duplicated "13/11" is corrected to
#    "14/11" (adds a new shared string).
# 3) Sprint 4!J6: a new task-time entry (1.5h) is logged for "Dev. Back-End"
#    on 14/11, which ripples into the shared B6 SUM formula.
# 4) Selection/active-cell bookkeeping: Tempo's cached selection moves to
#    N11, and Sprint 4 (the active/visible tab) ends with J7 selected.

$wb = $excel.ActiveWorkbook

$wsTempo = $wb.Worksheets.Item("Tempo")
$wsSprint4 = $wb.Worksheets.Item("Sprint 4")

# --- data fixes -----------------------------------------------------------

$wsTempo.Range("N4").Value = 50

$wsSprint4.Range("J2").Value = "14/11"
$wsSprint4.Range("J6").Value = 1.5

# --- selection bookkeeping (Tempo first, then land on Sprint 4 which stays
#     the active/visible tab, matching the saved workbook state) ----------

$wsTempo.Range("N11").Select()

$wsSprint4.Activate()
$wsSprint4.Range("J7").Select()
